# Auto-generated Excel COM-interop script applying the scheduled market-data/
# profit-recalculation refresh to the Seraph_Profits workbook.
#
# The workbook has 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), each a
# plain data table (Leve Name/Item/.../currentAveragePrice.../LeveProfit...) with
# no formulas - every cell holds a literal number pulled from an external price
# feed. This script rewrites the handful of rows whose market-price columns
# (H:N) were refreshed, updating each cell to its new value. A couple of cells
# that no longer have a value are cleared entirely.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 382.7143
$ws.Range("I11").Value = 382.7143
$ws.Range("K11").Value = 382.7143
$ws.Range("M11").Value = -242.7143
$ws.Range("H17").Value = 1666.1333
$ws.Range("I17").Value = 998.8
$ws.Range("J17").Value = 1999.8
$ws.Range("K17").Value = 2996.4
$ws.Range("L17").Value = 5999.4
$ws.Range("M17").Value = -2828.4
$ws.Range("N17").Value = -6335.4
$ws.Range("H40").Value = 2390
$ws.Range("J40").Value = 2975
$ws.Range("L40").Value = 2975
$ws.Range("N40").Value = -3325
$ws.Range("H41").Value = 375.0909
$ws.Range("J41").Value = 384
$ws.Range("L41").Value = 384
$ws.Range("N41").Value = -1264
$ws.Range("H51").Value = 4249.25
$ws.Range("I51").Value = 3998
$ws.Range("K51").Value = 3998
$ws.Range("M51").Value = -3514
$ws.Range("H55").Value = 519.8125
$ws.Range("J55").Value = 704.4545000000001
$ws.Range("L55").Value = 704.4545000000001
$ws.Range("N55").Value = -1132.4545
$ws.Range("H64").Value = 3199.75
$ws.Range("I64").Value = 3199
$ws.Range("K64").Value = 3199
$ws.Range("M64").Value = -2951
$ws.Range("H67").Value = 3199.75
$ws.Range("I67").Value = 3199
$ws.Range("K67").Value = 3199
$ws.Range("M67").Value = -2341
$ws.Range("H74").Value = 7051.5
$ws.Range("I74").Value = 2735.3333
$ws.Range("K74").Value = 2735.3333
$ws.Range("M74").Value = -1799.3333
$ws.Range("H76").Value = 6456.421
$ws.Range("J76").Value = 7090
$ws.Range("L76").Value = 7090
$ws.Range("N76").Value = -7720
$ws.Range("H77").Value = 7051.5
$ws.Range("I77").Value = 2735.3333
$ws.Range("K77").Value = 13676.6665
$ws.Range("M77").Value = -8996.666499999999
$ws.Range("H79").Value = 6456.421
$ws.Range("J79").Value = 7090
$ws.Range("L79").Value = 7090
$ws.Range("N79").Value = -9274
$ws.Range("H94").Value = 40000
$ws.Range("I94").Value = 40000
$ws.Range("K94").Value = 40000
$ws.Range("M94").Value = -39549
$ws.Range("H132").Value = 1760.5
$ws.Range("I132").Value = 1735.5
$ws.Range("K132").Value = 5206.5
$ws.Range("M132").Value = -2676.5
$ws.Range("H137").Value = 3008.2354
$ws.Range("I137").Value = 1436
$ws.Range("K137").Value = 4308
$ws.Range("M137").Value = -1758
$ws.Range("H138").Value = 5637.745
$ws.Range("I138").Value = 4221.625
$ws.Range("J138").Value = 5901.2095
$ws.Range("K138").Value = 12664.875
$ws.Range("L138").Value = 17703.6285
$ws.Range("M138").Value = -7524.875
$ws.Range("N138").Value = -27983.6285
$ws.Range("H141").Value = 2706.3333
$ws.Range("I141").Value = 2706.3333
$ws.Range("K141").Value = 8118.999899999999
$ws.Range("M141").Value = -2938.999899999999

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1805.5
$ws.Range("I61").Value = 1805.5
$ws.Range("K61").Value = 1805.5
$ws.Range("M61").Value = -1593.5
$ws.Range("H74").Value = 3533.375
$ws.Range("I74").Value = 976.75
$ws.Range("J74").Value = 6090
$ws.Range("K74").Value = 976.75
$ws.Range("L74").Value = 6090
$ws.Range("M74").Value = -102.75
$ws.Range("N74").Value = -7838
$ws.Range("H77").Value = 3533.375
$ws.Range("I77").Value = 976.75
$ws.Range("J77").Value = 6090
$ws.Range("K77").Value = 4883.75
$ws.Range("L77").Value = 30450
$ws.Range("M77").Value = -515.75
$ws.Range("N77").Value = -39186
$ws.Range("H97").Value = 999.6667
$ws.Range("I97").Value = 1274.75
$ws.Range("K97").Value = 1274.75
$ws.Range("M97").Value = -778.75
$ws.Range("H132").Value = 2018.9615
$ws.Range("I132").Value = 1862.5
$ws.Range("K132").Value = 5587.5
$ws.Range("M132").Value = -3057.5
$ws.Range("H136").Value = 1805.5
$ws.Range("I136").Value = 1805.5
$ws.Range("K136").Value = 5416.5
$ws.Range("M136").Value = -2866.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 3604.6667
$ws.Range("I25").Value = 3604.6667
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 3604.6667
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -3369.6667
$ws.Range("N25").ClearContents()
$ws.Range("H94").Value = 1500
$ws.Range("I94").Value = 1500
$ws.Range("K94").Value = 1500
$ws.Range("M94").Value = -1049
$ws.Range("H105").Value = 4538.852
$ws.Range("I105").Value = 4080.7693
$ws.Range("K105").Value = 4080.7693
$ws.Range("M105").Value = -2333.7693
$ws.Range("H107").Value = 1028.4
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 1455
$ws.Range("I134").Value = 699.375
$ws.Range("J134").Value = 7500
$ws.Range("K134").Value = 2098.125
$ws.Range("L134").Value = 22500
$ws.Range("M134").Value = 436.875
$ws.Range("N134").Value = -27570

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5846.4165
$ws.Range("I31").Value = 2574.5
$ws.Range("J31").Value = 7482.375
$ws.Range("K31").Value = 2574.5
$ws.Range("L31").Value = 7482.375
$ws.Range("M31").Value = -2279.5
$ws.Range("N31").Value = -8072.375
$ws.Range("H34").Value = 5846.4165
$ws.Range("I34").Value = 2574.5
$ws.Range("J34").Value = 7482.375
$ws.Range("K34").Value = 2574.5
$ws.Range("L34").Value = 7482.375
$ws.Range("M34").Value = -2372.5
$ws.Range("N34").Value = -7886.375
$ws.Range("H86").Value = 7904.5884
$ws.Range("J86").Value = 11600.375
$ws.Range("L86").Value = 11600.375
$ws.Range("N86").Value = -13846.375
$ws.Range("H89").Value = 7904.5884
$ws.Range("J89").Value = 11600.375
$ws.Range("L89").Value = 58001.875
$ws.Range("N89").Value = -69233.875
$ws.Range("H125").Value = 87374.75
$ws.Range("J125").Value = 87374.75
$ws.Range("L125").Value = 87374.75
$ws.Range("N125").Value = -92294.75
$ws.Range("H132").Value = 2883
$ws.Range("I132").Value = 2361.8235
$ws.Range("K132").Value = 7085.470499999999
$ws.Range("M132").Value = -4555.470499999999

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 35.25
$ws.Range("I10").Value = 35.25
$ws.Range("K10").Value = 105.75
$ws.Range("M10").Value = 33.25
$ws.Range("H34").Value = 1795.9333
$ws.Range("I34").Value = 1116.7778
$ws.Range("J34").Value = 2814.6667
$ws.Range("K34").Value = 3350.3334
$ws.Range("L34").Value = 8444.000100000001
$ws.Range("M34").Value = -3266.3334
$ws.Range("N34").Value = -8612.000100000001
$ws.Range("H113").Value = 956.8919
$ws.Range("I113").Value = 1531.1666
$ws.Range("K113").Value = 4593.4998
$ws.Range("M113").Value = -2423.4998

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1154.8611
$ws.Range("I102").Value = 508.37036
$ws.Range("K102").Value = 508.37036
$ws.Range("M102").Value = 1113.62964
$ws.Range("H113").Value = 3763.9285
$ws.Range("I113").Value = 2449.1667
$ws.Range("J113").Value = 4750
$ws.Range("K113").Value = 2449.1667
$ws.Range("L113").Value = 4750
$ws.Range("M113").Value = -279.1667000000002
$ws.Range("N113").Value = -9090
$ws.Range("H132").Value = 2096
$ws.Range("I132").Value = 2096
$ws.Range("K132").Value = 6288
$ws.Range("M132").Value = -3758

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3320
$ws.Range("I46").Value = 766.6667
$ws.Range("K46").Value = 766.6667
$ws.Range("M46").Value = -578.6667
$ws.Range("H55").Value = 355.69565
$ws.Range("J55").Value = 378.2857
$ws.Range("L55").Value = 378.2857
$ws.Range("N55").Value = -724.2857
$ws.Range("H68").Value = 2968.1875
$ws.Range("J68").Value = 2999.4
$ws.Range("L68").Value = 2999.4
$ws.Range("N68").Value = -4497.4
$ws.Range("H71").Value = 2968.1875
$ws.Range("J71").Value = 2999.4
$ws.Range("L71").Value = 14997
$ws.Range("N71").Value = -22485

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 12875
$ws.Range("I4").Value = 50000
$ws.Range("K4").Value = 50000
$ws.Range("M4").Value = -49887
$ws.Range("H46").Value = 244442.17
$ws.Range("J46").Value = 244442.17
$ws.Range("L46").Value = 244442.17
$ws.Range("N46").Value = -244904.17
$ws.Range("H96").Value = 1638.5
$ws.Range("I96").Value = 1582.75
$ws.Range("K96").Value = 1582.75
$ws.Range("M96").Value = -209.75
$ws.Range("H132").Value = 1440.9131
$ws.Range("I132").Value = 1156.1052
$ws.Range("K132").Value = 3468.3156
$ws.Range("M132").Value = -938.3155999999999
$ws.Range("H134").Value = 244442.17
$ws.Range("J134").Value = 244442.17
$ws.Range("L134").Value = 733326.51
$ws.Range("N134").Value = -738396.51

